$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 data values: replace the sample/placeholder data with the real
# (Arabic) values, keep just plain text where applicable.
$ws.Range("A2").Value = "عاطف العباسي"
$ws.Range("B2").Value = "no5510425@gmail.com"
$ws.Range("C2").Value = "مهندس"
$ws.Range("D2").Value = "الاستاذ"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = "crt_Part_1.pdf"
$ws.Range("G2").Value = "Abbasi's-affection"

# B2 now holds an e-mail address - turn it into a real mailto hyperlink
# (this also creates the built-in "Hyperlink" style/font Excel uses).
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:no5510425@gmail.com")

# Update the active selection to match the new focus cell.
$ws.Range("B2").Select()

Write-Output "done"
